$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = -12.1309
$ws.Range("C7").Value = -12.8387
$ws.Range("D7").Value = -7.327899999999995
$ws.Range("D15").Value = -8.361399999999994
$ws.Range("C16").Value = -14.0645
$ws.Range("D21").Value = -8.533399999999995
$ws.Range("D22").Value = -8.066400000000007
$ws.Range("D23").Value = -7.347299999999998
$ws.Range("C28").Value = -12.5759
$ws.Range("C29").Value = -11.30220000000001
$ws.Range("C32").Value = -13.22940000000001
$ws.Range("D34").Value = -7.989199999999999
$ws.Range("C40").Value = -12.88620000000001
$ws.Range("D43").Value = -8.266700000000004
$ws.Range("D45").Value = -7.794599999999996
$ws.Range("D50").Value = -8.167399999999997
$ws.Range("D51").Value = -7.532499999999999
$ws.Range("C52").Value = -11.22130000000001
$ws.Range("C57").Value = -13.85839999999999
$ws.Range("C66").Value = -12.1791
$ws.Range("D66").Value = -7.9534
$ws.Range("D67").Value = -6.3423
$ws.Range("D79").Value = -6.2864
$ws.Range("D84").Value = -8.627200000000004
$ws.Range("D92").Value = -6.433400000000002
$ws.Range("D97").Value = -8.271199999999999
$ws.Range("C100").Value = -12.4551
